# 678-MS-EPP-DB-DL-REC-NON-RNI-CTPD-DL-MD-TR-1-EarlyRePayment-Newcreateloan.xlsx
# "Loan RBI, Variable Instalments"
#
# On the "Repayment Schedule" sheet, insert a new (blank) column before the
# existing "Late" column so the schedule gains a spare "Variable
# Instalments" column ahead of Late / Outstanding, then leave that sheet
# as the active tab (matching the workbook's new activeTab/tabSelected
# state).

$wb = $excel.ActiveWorkbook

$wsRepayment = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new column N (pushing the old "Late" column from N to O, and
# "Outstanding" from O to P/Q accordingly).
$wsRepayment.Columns("N:N").Insert()

# Match the new column's width to its left neighbour (column M).
$wsRepayment.Columns("N:N").ColumnWidth = 10.33

# Make "Repayment Schedule" the active (selected) tab, as it was in the
# edited workbook, and update the active cell/selection on it.
$wsRepayment.Activate() | Out-Null
$wsRepayment.Range("T9").Select() | Out-Null
